$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c0 = $ws.Range("D2")
$c0.NumberFormat = "@"
$c0.Value = "26.827.12"
$c0.ClearFormats()

$c1 = $ws.Range("E2")
$c1.NumberFormat = "@"
$c1.Value = "  +0.07%  "
$c1.ClearFormats()

$c2 = $ws.Range("D3")
$c2.NumberFormat = "@"
$c2.Value = "1.542.74"
$c2.ClearFormats()

$c3 = $ws.Range("E3")
$c3.NumberFormat = "@"
$c3.Value = "  -1.57%  "
$c3.ClearFormats()

$c4 = $ws.Range("E4")
$c4.NumberFormat = "@"
$c4.Value = "  +0.20%  "
$c4.ClearFormats()

$c5 = $ws.Range("D5")
$c5.NumberFormat = "@"
$c5.Value = "205.94"
$c5.ClearFormats()

$c6 = $ws.Range("E5")
$c6.NumberFormat = "@"
$c6.Value = "  -0.28%  "
$c6.ClearFormats()

$c7 = $ws.Range("E6")
$c7.NumberFormat = "@"
$c7.Value = "  -0.63%  "
$c7.ClearFormats()

$c8 = $ws.Range("E7")
$c8.NumberFormat = "@"
$c8.Value = "  +0.18%  "
$c8.ClearFormats()

$c9 = $ws.Range("E8")
$c9.NumberFormat = "@"
$c9.Value = "  -0.52%  "
$c9.ClearFormats()

$c10 = $ws.Range("D9")
$c10.NumberFormat = "@"
$c10.Value = "21.42"
$c10.ClearFormats()

$c11 = $ws.Range("E9")
$c11.NumberFormat = "@"
$c11.Value = "  -2.67%  "
$c11.ClearFormats()

$c12 = $ws.Range("D10")
$c12.NumberFormat = "@"
$c12.Value = "0.0582"
$c12.ClearFormats()

$c13 = $ws.Range("E10")
$c13.NumberFormat = "@"
$c13.Value = "  -0.48%  "
$c13.ClearFormats()

$c14 = $ws.Range("E11")
$c14.NumberFormat = "@"
$c14.Value = "  -1.05%  "
$c14.ClearFormats()

$c15 = $ws.Range("D12")
$c15.NumberFormat = "@"
$c15.Value = "1.761.08"
$c15.ClearFormats()

$c16 = $ws.Range("E12")
$c16.NumberFormat = "@"
$c16.Value = "  -1.62%  "
$c16.ClearFormats()

$c17 = $ws.Range("D13")
$c17.NumberFormat = "@"
$c17.Value = "1.541.15"
$c17.ClearFormats()

$c18 = $ws.Range("E13")
$c18.NumberFormat = "@"
$c18.Value = "  -1.61%  "
$c18.ClearFormats()

$c19 = $ws.Range("E14")
$c19.NumberFormat = "@"
$c19.Value = "  -1.46%  "
$c19.ClearFormats()

$c20 = $ws.Range("E15")
$c20.NumberFormat = "@"
$c20.Value = "  -0.79%  "
$c20.ClearFormats()

$c21 = $ws.Range("D16")
$c21.NumberFormat = "@"
$c21.Value = "26.821.07"
$c21.ClearFormats()

$c22 = $ws.Range("E16")
$c22.NumberFormat = "@"
$c22.Value = "  +0.05%  "
$c22.ClearFormats()

$c23 = $ws.Range("D17")
$c23.NumberFormat = "@"
$c23.Value = "61.27"
$c23.ClearFormats()

$c24 = $ws.Range("E17")
$c24.NumberFormat = "@"
$c24.Value = "  -0.27%  "
$c24.ClearFormats()

$c25 = $ws.Range("D18")
$c25.NumberFormat = "@"
$c25.Value = "214.77"
$c25.ClearFormats()

$c26 = $ws.Range("D19")
$c26.NumberFormat = "@"
$c26.Value = "7.25"
$c26.ClearFormats()

$c27 = $ws.Range("E19")
$c27.NumberFormat = "@"
$c27.Value = "  -2.28%  "
$c27.ClearFormats()

$c28 = $ws.Range("E20")
$c28.NumberFormat = "@"
$c28.Value = "  +0.69%  "
$c28.ClearFormats()

$c29 = $ws.Range("D22")
$c29.NumberFormat = "@"
$c29.Value = "4.00"
$c29.ClearFormats()

$c30 = $ws.Range("E22")
$c30.NumberFormat = "@"
$c30.Value = "  -3.04%  "
$c30.ClearFormats()

$c31 = $ws.Range("E23")
$c31.NumberFormat = "@"
$c31.Value = "  -1.29%  "
$c31.ClearFormats()

$c32 = $ws.Range("D24")
$c32.NumberFormat = "@"
$c32.Value = "1.94"
$c32.ClearFormats()

$c33 = $ws.Range("E24")
$c33.NumberFormat = "@"
$c33.Value = "  -3.01%  "
$c33.ClearFormats()

$c34 = $ws.Range("D25")
$c34.NumberFormat = "@"
$c34.Value = "152.78"
$c34.ClearFormats()

$c35 = $ws.Range("E25")
$c35.NumberFormat = "@"
$c35.Value = "  -0.40%  "
$c35.ClearFormats()

$c36 = $ws.Range("E26")
$c36.NumberFormat = "@"
$c36.Value = "  -2.25%  "
$c36.ClearFormats()

$c37 = $ws.Range("D27")
$c37.NumberFormat = "@"
$c37.Value = "14.83"
$c37.ClearFormats()

$c38 = $ws.Range("E27")
$c38.NumberFormat = "@"
$c38.Value = "  -0.89%  "
$c38.ClearFormats()

$c39 = $ws.Range("E28")
$c39.NumberFormat = "@"
$c39.Value = "  +0.14%  "
$c39.ClearFormats()

$c40 = $ws.Range("E29")
$c40.NumberFormat = "@"
$c40.Value = "  -0.86%  "
$c40.ClearFormats()

$c41 = $ws.Range("B30")
$c41.Value = "PancakeSwap"

$c42 = $ws.Range("C30")
$c42.Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"

$c43 = $ws.Range("D30")
$c43.NumberFormat = "@"
$c43.Value = "1.10"
$c43.ClearFormats()

$c44 = $ws.Range("E30")
$c44.NumberFormat = "@"
$c44.Value = "  -1.24%  "
$c44.ClearFormats()

$c45 = $ws.Range("B31")
$c45.Value = "Hedera"

$c46 = $ws.Range("C31")
$c46.Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"

$c47 = $ws.Range("D31")
$c47.NumberFormat = "@"
$c47.Value = "0.0458"
$c47.ClearFormats()

$c48 = $ws.Range("E31")
$c48.NumberFormat = "@"
$c48.Value = "  -1.83%  "
$c48.ClearFormats()

$c49 = $ws.Range("E32")
$c49.NumberFormat = "@"
$c49.Value = "  +1.68%  "
$c49.ClearFormats()

$c50 = $ws.Range("D33")
$c50.NumberFormat = "@"
$c50.Value = "1.368.56"
$c50.ClearFormats()

$c51 = $ws.Range("E33")
$c51.NumberFormat = "@"
$c51.Value = "  -2.09%  "
$c51.ClearFormats()

$c52 = $ws.Range("E34")
$c52.NumberFormat = "@"
$c52.Value = "  +0.34%  "
$c52.ClearFormats()

$c53 = $ws.Range("D35")
$c53.NumberFormat = "@"
$c53.Value = "1.51"
$c53.ClearFormats()

$c54 = $ws.Range("E35")
$c54.NumberFormat = "@"
$c54.Value = "  -1.29%  "
$c54.ClearFormats()

$c55 = $ws.Range("D36")
$c55.NumberFormat = "@"
$c55.Value = "0.965"
$c55.ClearFormats()

$c56 = $ws.Range("E36")
$c56.NumberFormat = "@"
$c56.Value = "  +3.13%  "
$c56.ClearFormats()

$c57 = $ws.Range("E37")
$c57.NumberFormat = "@"
$c57.Value = "  -0.01%  "
$c57.ClearFormats()

$c58 = $ws.Range("E38")
$c58.NumberFormat = "@"
$c58.Value = "  +0.79%  "
$c58.ClearFormats()

$c59 = $ws.Range("D39")
$c59.NumberFormat = "@"
$c59.Value = "0.521"
$c59.ClearFormats()

$c60 = $ws.Range("E39")
$c60.NumberFormat = "@"
$c60.Value = "  -1.56%  "
$c60.ClearFormats()

$c61 = $ws.Range("E40")
$c61.NumberFormat = "@"
$c61.Value = "  +8.55%  "
$c61.ClearFormats()

$c62 = $ws.Range("D41")
$c62.NumberFormat = "@"
$c62.Value = "0.806"
$c62.ClearFormats()

$c63 = $ws.Range("E41")
$c63.NumberFormat = "@"
$c63.Value = "  -1.03%  "
$c63.ClearFormats()

$c64 = $ws.Range("B42")
$c64.Value = "PaxDollar"

$c65 = $ws.Range("C42")
$c65.Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"

$c66 = $ws.Range("D42")
$c66.NumberFormat = "@"
$c66.Value = "1.00"
$c66.ClearFormats()

$c67 = $ws.Range("E42")
$c67.NumberFormat = "@"
$c67.Value = "  +0.16%  "
$c67.ClearFormats()

$c68 = $ws.Range("B43")
$c68.Value = "WEMIXToken"

$c69 = $ws.Range("C43")
$c69.Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"

$c70 = $ws.Range("D43")
$c70.NumberFormat = "@"
$c70.Value = "0.991"
$c70.ClearFormats()

$c71 = $ws.Range("E43")
$c71.NumberFormat = "@"
$c71.Value = "  +0.35%  "
$c71.ClearFormats()

$c72 = $ws.Range("B44")
$c72.Value = "MXToken"

$c73 = $ws.Range("C44")
$c73.Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"

$c74 = $ws.Range("D44")
$c74.NumberFormat = "@"
$c74.Value = "2.21"
$c74.ClearFormats()

$c75 = $ws.Range("E44")
$c75.NumberFormat = "@"
$c75.Value = "  +1.01%  "
$c75.ClearFormats()

$c76 = $ws.Range("B45")
$c76.Value = "Aave"

$c77 = $ws.Range("C45")
$c77.Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"

$c78 = $ws.Range("D45")
$c78.NumberFormat = "@"
$c78.Value = "63.16"
$c78.ClearFormats()

$c79 = $ws.Range("E45")
$c79.NumberFormat = "@"
$c79.Value = "  -0.27%  "
$c79.ClearFormats()

$c80 = $ws.Range("B46")
$c80.Value = "RenderToken"

$c81 = $ws.Range("C46")
$c81.Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"

$c82 = $ws.Range("D46")
$c82.NumberFormat = "@"
$c82.Value = "1.74"
$c82.ClearFormats()

$c83 = $ws.Range("E46")
$c83.NumberFormat = "@"
$c83.Value = "  -3.46%  "
$c83.ClearFormats()

$c84 = $ws.Range("B47")
$c84.Value = "RocketPoolETH"

$c85 = $ws.Range("C47")
$c85.Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"

$c86 = $ws.Range("D47")
$c86.NumberFormat = "@"
$c86.Value = "1.675.59"
$c86.ClearFormats()

$c87 = $ws.Range("E47")
$c87.NumberFormat = "@"
$c87.Value = "  -1.65%  "
$c87.ClearFormats()

$c88 = $ws.Range("B48")
$c88.Value = "Quant"

$c89 = $ws.Range("C48")
$c89.Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"

$c90 = $ws.Range("D48")
$c90.NumberFormat = "@"
$c90.Value = "84.40"
$c90.ClearFormats()

$c91 = $ws.Range("E48")
$c91.NumberFormat = "@"
$c91.Value = "  -1.88%  "
$c91.ClearFormats()

$c92 = $ws.Range("B49")
$c92.Value = "Cronos"

$c93 = $ws.Range("C49")
$c93.Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"

$c94 = $ws.Range("D49")
$c94.NumberFormat = "@"
$c94.Value = "0.0509"
$c94.ClearFormats()

$c95 = $ws.Range("E49")
$c95.NumberFormat = "@"
$c95.Value = "  +3.50%  "
$c95.ClearFormats()

$c96 = $ws.Range("B50")
$c96.Value = "BabyDogeCoin"

$c97 = $ws.Range("C50")
$c97.Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"

$c98 = $ws.Range("D50")
$c98.NumberFormat = "@"
$c98.Value = "0.0₇0977"
$c98.ClearFormats()

$c99 = $ws.Range("E50")
$c99.NumberFormat = "@"
$c99.Value = "  -0.70%  "
$c99.ClearFormats()

$c100 = $ws.Range("D51")
$c100.NumberFormat = "@"
$c100.Value = "0.0941"
$c100.ClearFormats()

$c101 = $ws.Range("E51")
$c101.NumberFormat = "@"
$c101.Value = "  -1.42%  "
$c101.ClearFormats()

